$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated query / filename text blocks (filter values changed to Medullary Carcinoma /
# 16-20 / (2,3] / tumor_grade Not Reported; filenames renamed TC02->TC03 TumorGrade)
$tc03Neo4j = @'
TC03_Bento_MultiFilter_Diagnosis-Recurrence-TumorSize-TumorGrade_Neo4jData.xlsx
'@
$tc03Web = @'
TC03_Bento_MultiFilter_Diagnosis-Recurrence-TumorSize-TumorGrade_WebData.xlsx
'@
$countQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"]and d.tumor_size_group In ["(2,3]"] and d.tumor_grade In ["Not Reported"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
'@
$filesQuery = @'
MATCH (f:file)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"]and d.tumor_size_group In ["(2,3]"] and d.tumor_grade In ["Not Reported"]
RETURN  f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
    f.file_size AS `Size`,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`,
    ss.disease_subtype as `Diagnosis`
'@
$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"]and d.tumor_size_group In ["(2,3]"] and d.tumor_grade In ["Not Reported"]
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`,
            lp.test_name as Platform
'@
$casesQuery = @'
MATCH (ss:study_subject)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, collect(DISTINCT f) AS files
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
 WHERE ss.disease_subtype IN ["Medullary Carcinoma"] and sf.grouped_recurrence_score IN ["16-20"]and d.tumor_size_group In ["(2,3]"] and d.tumor_grade In ["Not Reported"]
return ss.study_subject_id as `Case ID`,
       p.program_acronym as `Program Code`,
        p.program_id as Program_ID,
       s.study_acronym as `Arm`,
       ss.disease_subtype as `Diagnosis`,
       sf.grouped_recurrence_score AS `Recurrence Score`,
       d.tumor_size_group AS `tumor_size`,
       d.er_status AS `ER Status`,
       d.pr_status AS `PR Status`,
       demo.age_at_index AS `Age (years)`,
demo.survival_time AS `Survival (days)`
'@

# Row 2 - CasesTab
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $countQuery
$ws.Range("D2").Value = $tc03Neo4j
$ws.Range("E2").Value = $tc03Web

# Row 3 - SamplesTab
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $countQuery
$ws.Range("D3").Value = $tc03Neo4j
$ws.Range("E3").Value = $tc03Web

# Row 4 - FilesTab
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $countQuery
$ws.Range("D4").Value = $tc03Neo4j
$ws.Range("E4").Value = $tc03Web

# Row 2 grew taller to fit the longer wrapped Cases query text
$ws.Range("B2:E2").RowHeight = 348

# Restore the view that was saved with the workbook: scrolled so row 4 is the
# top visible row, with C4 as the active selected cell
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C4").Select() | Out-Null
